# Updates the "cryptos" price-ranking sheet (Sheet1) to the refreshed
# snapshot: new Price (column D) / Volume(1h) (column E) readings for most
# rows, plus the coin in rank position 43/44 (rows 44/45) swapping places
# (Hedera <-> InjectiveProtocol, including its link + price + change%).
#
# Price values that still *look* numeric (e.g. "583.04") would otherwise be
# silently reinterpreted by Excel as a Double when assigned through
# .Value, losing the original text formatting/precision that this sheet
# relies on (prices like "60.911.70" use dots as thousand separators, not
# decimals). Set-TextValue guards against that by forcing the cell to a
# text format before the write, then restores the default "Normal" style
# so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '60.911.70'
$ws.Range('E2').Value = '  +0.64%  '

$ws.Range('D3').Value = '2.639.55'
$ws.Range('E3').Value = '  +1.45%  '

$ws.Range('E4').Value = '  -0.11%  '

Set-TextValue $ws.Range('D5') '583.04'
$ws.Range('E5').Value = '  -0.09%  '

Set-TextValue $ws.Range('D6') '144.02'
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('E8').Value = '  -0.37%  '

$ws.Range('E9').Value = '  +0.52%  '

$ws.Range('E10').Value = '  +2.08%  '

$ws.Range('E11').Value = '  +1.27%  '

$ws.Range('E12').Value = '  +1.18%  '

$ws.Range('D13').Value = '3.106.55'
$ws.Range('E13').Value = '  +0.69%  '

Set-TextValue $ws.Range('D14') '26.23'
$ws.Range('E14').Value = '  +6.84%  '

$ws.Range('D15').Value = '60.867.79'
$ws.Range('E15').Value = '  +0.57%  '

$ws.Range('E16').Value = '  +1.70%  '

$ws.Range('D17').Value = '2.653.11'
$ws.Range('E17').Value = '  +1.09%  '

$ws.Range('E18').Value = '  +2.04%  '

$ws.Range('E19').Value = '  +0.99%  '

Set-TextValue $ws.Range('D20') '351.73'
$ws.Range('E20').Value = '  +0.95%  '

Set-TextValue $ws.Range('D21') '6.87'
$ws.Range('E21').Value = '  -0.58%  '

Set-TextValue $ws.Range('D22') '1.00'
$ws.Range('E22').Value = '  +0.03%  '

Set-TextValue $ws.Range('D23') '0.525'
$ws.Range('E23').Value = '  +1.01%  '

Set-TextValue $ws.Range('D24') '63.95'
$ws.Range('E24').Value = '  +1.34%  '

$ws.Range('E25').Value = '  +1.51%  '

$ws.Range('E26').Value = '  -0.55%  '

Set-TextValue $ws.Range('D27') '8.37'
$ws.Range('E27').Value = '  +5.42%  '

$ws.Range('E28').Value = '  +6.65%  '

$ws.Range('D29').Value = '0.0₃0809'
$ws.Range('E29').Value = '  +1.39%  '

Set-TextValue $ws.Range('D30') '6.75'
$ws.Range('E30').Value = '  +5.77%  '

Set-TextValue $ws.Range('D31') '167.34'
$ws.Range('E31').Value = '  +1.74%  '

$ws.Range('E32').Value = '  -0.05%  '

Set-TextValue $ws.Range('D33') '19.96'
$ws.Range('E33').Value = '  +2.32%  '

Set-TextValue $ws.Range('D34') '4.59'
$ws.Range('E34').Value = '  +7.87%  '

Set-TextValue $ws.Range('D35') '1.08'
$ws.Range('E35').Value = '  +7.97%  '

$ws.Range('E36').Value = '  +7.15%  '

$ws.Range('E37').Value = '  +3.95%  '

Set-TextValue $ws.Range('D38') '342.34'
$ws.Range('E38').Value = '  +9.24%  '

Set-TextValue $ws.Range('D39') '4.11'
$ws.Range('E39').Value = '  +5.70%  '

Set-TextValue $ws.Range('D40') '0.902'
$ws.Range('E40').Value = '  +7.09%  '

Set-TextValue $ws.Range('D41') '38.26'
$ws.Range('E41').Value = '  +0.79%  '

Set-TextValue $ws.Range('D42') '138.67'
$ws.Range('E42').Value = '  +2.48%  '

$ws.Range('E43').Value = '  +3.77%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D44') '21.11'
$ws.Range('E44').Value = '  +4.38%  '

$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D45') '0.0573'
$ws.Range('E45').Value = '  +3.66%  '

Set-TextValue $ws.Range('D46') '0.626'
$ws.Range('E46').Value = '  +3.05%  '

Set-TextValue $ws.Range('D47') '20.27'
$ws.Range('E47').Value = '  +2.00%  '

$ws.Range('E48').Value = '  +2.98%  '

$ws.Range('E49').Value = '  +0.49%  '

Set-TextValue $ws.Range('D50') '0.998'
$ws.Range('E50').Value = '  +0.18%  '

$ws.Range('D51').Value = '2.089.13'
$ws.Range('E51').Value = '  +2.57%  '
